$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '66.221.18'
$ws.Range("E2").Value = '  -3.83%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.548.61'
$ws.Range("E3").Value = '  -4.54%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
Set-TextValue $ws.Range("D5") '585.37'
$ws.Range("E5").Value = '  -3.24%  '

# Row 6
Set-TextValue $ws.Range("D6") '180.28'
$ws.Range("E6").Value = '  -1.38%  '

# Row 7
Set-TextValue $ws.Range("D7") '3.541.04'
$ws.Range("E7").Value = '  -4.68%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.604'
$ws.Range("E8").Value = '  -4.38%  '

# Row 9
$ws.Range("E9").Value = '  +0.46%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.662'
$ws.Range("E10").Value = '  -8.04%  '

# Row 11
Set-TextValue $ws.Range("D11") '53.01'
$ws.Range("E11").Value = '  -6.91%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.141'
$ws.Range("E12").Value = '  -12.41%  '

# Row 13
Set-TextValue $ws.Range("D13") '0.0000247'
$ws.Range("E13").Value = '  -16.25%  '

# Row 14
Set-TextValue $ws.Range("D14") '9.68'
$ws.Range("E14").Value = '  -9.51%  '

# Row 15
Set-TextValue $ws.Range("D15") '4.134.47'
$ws.Range("E15").Value = '  -4.23%  '

# Row 16
Set-TextValue $ws.Range("D16") '3.557.53'
$ws.Range("E16").Value = '  -4.18%  '

# Row 17
$ws.Range("E17").Value = '  -0.81%  '

# Row 18
Set-TextValue $ws.Range("D18") '18.20'
$ws.Range("E18").Value = '  -6.70%  '

# Row 19
Set-TextValue $ws.Range("D19") '66.126.39'
$ws.Range("E19").Value = '  -3.95%  '

# Row 20
Set-TextValue $ws.Range("D20") '12.02'
$ws.Range("E20").Value = '  -7.46%  '

# Row 21
Set-TextValue $ws.Range("D21") '1.04'
$ws.Range("E21").Value = '  -7.74%  '

# Row 22
Set-TextValue $ws.Range("D22") '389.60'
$ws.Range("E22").Value = '  -6.12%  '

# Row 23
Set-TextValue $ws.Range("D23") '4.26'
$ws.Range("E23").Value = '  -8.94%  '

# Row 24
Set-TextValue $ws.Range("D24") '84.54'
$ws.Range("E24").Value = '  -5.19%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.84'
$ws.Range("E25").Value = '  -7.27%  '

# Row 26
Set-TextValue $ws.Range("D26") '12.18'
$ws.Range("E26").Value = '  -4.45%  '

# Row 27
$ws.Range("E27").Value = '  -1.03%  '

# Row 28
Set-TextValue $ws.Range("D28") '10.12'
$ws.Range("E28").Value = '  -7.75%  '

# Row 29
Set-TextValue $ws.Range("D29") '3.54'
$ws.Range("E29").Value = '  -9.05%  '

# Row 30
Set-TextValue $ws.Range("D30") '8.82'
$ws.Range("E30").Value = '  -8.53%  '

# Row 31
Set-TextValue $ws.Range("D31") '30.72'
$ws.Range("E31").Value = '  -7.19%  '

# Row 32
Set-TextValue $ws.Range("D32") '6.65'
$ws.Range("E32").Value = '  -9.47%  '

# Row 33
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D33") '65.08'
$ws.Range("E33").Value = '  +0.12%  '

# Row 34
Set-TextValue $ws.Range("D34") '11.86'
$ws.Range("E34").Value = '  -5.17%  '

# Row 35
$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D35") '607.00'
$ws.Range("E35").Value = '  +0.24%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.111'
$ws.Range("E36").Value = '  -7.78%  '

# Row 37
Set-TextValue $ws.Range("D37") '40.91'
$ws.Range("E37").Value = '  -7.05%  '

# Row 38
Set-TextValue $ws.Range("D38") '1.00'
$ws.Range("E38").Value = '  +0.00%  '

# Row 39
$ws.Range("E39").Value = '  -0.09%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.368'
$ws.Range("E40").Value = '  -9.72%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.0₃0732'
$ws.Range("E41").Value = '  -17.49%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.129'
$ws.Range("E42").Value = '  -7.31%  '

# Row 43
Set-TextValue $ws.Range("D43") '2.858.12'
$ws.Range("E43").Value = '  +2.92%  '

# Row 44
Set-TextValue $ws.Range("D44") '2.76'
$ws.Range("E44").Value = '  -10.40%  '

# Row 45
Set-TextValue $ws.Range("D45") '0.0403'
$ws.Range("E45").Value = '  -9.02%  '

# Row 46
Set-TextValue $ws.Range("D46") '2.38'
$ws.Range("E46").Value = '  -10.72%  '

# Row 47
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D47") '3.08'
$ws.Range("E47").Value = '  -0.31%  '

# Row 48
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D48") '0.129'
$ws.Range("E48").Value = '  -5.01%  '

# Row 49
Set-TextValue $ws.Range("D49") '136.80'
$ws.Range("E49").Value = '  -2.87%  '

# Row 50
Set-TextValue $ws.Range("D50") '2.46'
$ws.Range("E50").Value = '  -10.02%  '

# Row 51
Set-TextValue $ws.Range("D51") '8.18'
$ws.Range("E51").Value = '  -11.55%  '
